$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A2").Value = "Última actualización: 08:39:56"
$ws.Range("A3").Value = "Total filas: 118"
$data = New-Object 'object[,]' 118,5
$data[0,0] = "04:01:01"
$data[0,1] = "04:01"
$data[0,2] = "81_EL PELIGRO"
$data[0,3] = 0
$data[0,4] = "LP1912"
$data[1,0] = "04:36:47"
$data[1,1] = "04:45"
$data[1,2] = "215A_EL PATO"
$data[1,3] = 9
$data[1,4] = "LP1912"
$data[2,0] = "04:01:01"
$data[2,1] = "04:46"
$data[2,2] = "215A_EL PATO"
$data[2,3] = 45
$data[2,4] = "LP1912"
$data[3,0] = "04:01:01"
$data[3,1] = "04:53"
$data[3,2] = "11_ETCHEVERRY"
$data[3,3] = 52
$data[3,4] = "LP1912"
$data[4,0] = "04:51:28"
$data[4,1] = "05:13"
$data[4,2] = "14_ABASTO"
$data[4,3] = 22
$data[4,4] = "LP1912"
$data[5,0] = "04:36:47"
$data[5,1] = "05:14"
$data[5,2] = "14_ABASTO"
$data[5,3] = 38
$data[5,4] = "LP1912"
$data[6,0] = "04:01:01"
$data[6,1] = "05:16"
$data[6,2] = "17_ROMERO"
$data[6,3] = 75
$data[6,4] = "LP1912"
$data[7,0] = "05:20:00"
$data[7,1] = "05:22"
$data[7,2] = "14_ABASTO"
$data[7,3] = 2
$data[7,4] = "LP1912"
$data[8,0] = "04:01:01"
$data[8,1] = "05:22"
$data[8,2] = "23_HERNANDEZ"
$data[8,3] = 81
$data[8,4] = "LP1912"
$data[9,0] = "04:36:47"
$data[9,1] = "05:34"
$data[9,2] = "215B_EL PATO"
$data[9,3] = 58
$data[9,4] = "LP1912"
$data[10,0] = "04:01:01"
$data[10,1] = "05:35"
$data[10,2] = "215B_EL PATO"
$data[10,3] = 94
$data[10,4] = "LP1912"
$data[11,0] = "04:01:01"
$data[11,1] = "05:41"
$data[11,2] = "14_ABASTO"
$data[11,3] = 100
$data[11,4] = "LP1912"
$data[12,0] = "04:01:01"
$data[12,1] = "05:46"
$data[12,2] = "15_ABASTO"
$data[12,3] = 105
$data[12,4] = "LP1912"
$data[13,0] = "04:36:47"
$data[13,1] = "06:04"
$data[13,2] = "16_SANTA ANA"
$data[13,3] = 88
$data[13,4] = "LP1912"
$data[14,0] = "04:36:47"
$data[14,1] = "06:11"
$data[14,2] = "215A_EL PATO"
$data[14,3] = 95
$data[14,4] = "LP1912"
$data[15,0] = "05:20:00"
$data[15,1] = "06:12"
$data[15,2] = "215A_EL PATO"
$data[15,3] = 52
$data[15,4] = "LP1912"
$data[16,0] = "04:36:47"
$data[16,1] = "06:14"
$data[16,2] = "225_HARAS DEL SUR"
$data[16,3] = 98
$data[16,4] = "LP1912"
$data[17,0] = "04:36:47"
$data[17,1] = "06:21"
$data[17,2] = "26_HERNANDEZ"
$data[17,3] = 105
$data[17,4] = "LP1912"
$data[18,0] = "04:36:47"
$data[18,1] = "06:27"
$data[18,2] = "23_HERNANDEZ"
$data[18,3] = 111
$data[18,4] = "LP1912"
$data[19,0] = "04:36:47"
$data[19,1] = "06:29"
$data[19,2] = "86_EST CHICA-ESC AGRARIA"
$data[19,3] = 113
$data[19,4] = "LP1912"
$data[20,0] = "05:20:00"
$data[20,1] = "06:30"
$data[20,2] = "86_EST CHICA-ESC AGRARIA"
$data[20,3] = 70
$data[20,4] = "LP1912"
$data[21,0] = "04:36:47"
$data[21,1] = "06:31"
$data[21,2] = "16_SANTA ANA"
$data[21,3] = 115
$data[21,4] = "LP1912"
$data[22,0] = "06:23:52"
$data[22,1] = "06:33"
$data[22,2] = "23_HERNANDEZ"
$data[22,3] = 10
$data[22,4] = "LP1912"
$data[23,0] = "06:23:52"
$data[23,1] = "06:43"
$data[23,2] = "225_C ROCA-H SUR"
$data[23,3] = 20
$data[23,4] = "LP1912"
$data[24,0] = "04:51:28"
$data[24,1] = "06:44"
$data[24,2] = "225_C ROCA-H SUR"
$data[24,3] = 113
$data[24,4] = "LP1912"
$data[25,0] = "04:51:28"
$data[25,1] = "06:46"
$data[25,2] = "215C_EL PATO"
$data[25,3] = 115
$data[25,4] = "LP1912"
$data[26,0] = "05:20:00"
$data[26,1] = "06:47"
$data[26,2] = "215C_EL PATO"
$data[26,3] = 87
$data[26,4] = "LP1912"
$data[27,0] = "06:52:23"
$data[27,1] = "06:52"
$data[27,2] = "215C_EL PATO"
$data[27,3] = 0
$data[27,4] = "LP1912"
$data[28,0] = "05:51:32"
$data[28,1] = "06:59"
$data[28,2] = "14_ABASTO"
$data[28,3] = 68
$data[28,4] = "LP1912"
$data[29,0] = "05:20:00"
$data[29,1] = "07:00"
$data[29,2] = "10_OLMOS"
$data[29,3] = 100
$data[29,4] = "LP1912"
$data[30,0] = "05:20:00"
$data[30,1] = "07:00"
$data[30,2] = "14_ABASTO"
$data[30,3] = 100
$data[30,4] = "LP1912"
$data[31,0] = "06:23:52"
$data[31,1] = "07:01"
$data[31,2] = "16_SANTA ANA"
$data[31,3] = 38
$data[31,4] = "LP1912"
$data[32,0] = "05:51:32"
$data[32,1] = "07:04"
$data[32,2] = "23_HERNANDEZ"
$data[32,3] = 73
$data[32,4] = "LP1912"
$data[33,0] = "05:20:00"
$data[33,1] = "07:05"
$data[33,2] = "15_ABASTO"
$data[33,3] = 105
$data[33,4] = "LP1912"
$data[34,0] = "06:23:52"
$data[34,1] = "07:06"
$data[34,2] = "225_GOMEZ"
$data[34,3] = 43
$data[34,4] = "LP1912"
$data[35,0] = "05:20:00"
$data[35,1] = "07:07"
$data[35,2] = "225_GOMEZ"
$data[35,3] = 107
$data[35,4] = "LP1912"
$data[36,0] = "05:51:32"
$data[36,1] = "07:11"
$data[36,2] = "215A_EL PATO"
$data[36,3] = 80
$data[36,4] = "LP1912"
$data[37,0] = "05:20:00"
$data[37,1] = "07:12"
$data[37,2] = "215A_EL PATO"
$data[37,3] = 112
$data[37,4] = "LP1912"
$data[38,0] = "05:51:32"
$data[38,1] = "07:15"
$data[38,2] = "11_ETCHEVERRY"
$data[38,3] = 84
$data[38,4] = "LP1912"
$data[39,0] = "06:52:23"
$data[39,1] = "07:16"
$data[39,2] = "16_SANTA ANA"
$data[39,3] = 24
$data[39,4] = "LP1912"
$data[40,0] = "05:20:00"
$data[40,1] = "07:16"
$data[40,2] = "11_ETCHEVERRY"
$data[40,3] = 116
$data[40,4] = "LP1912"
$data[41,0] = "05:51:32"
$data[41,1] = "07:21"
$data[41,2] = "26_HERNANDEZ"
$data[41,3] = 90
$data[41,4] = "LP1912"
$data[42,0] = "06:23:52"
$data[42,1] = "07:22"
$data[42,2] = "10_OLMOS"
$data[42,3] = 59
$data[42,4] = "LP1912"
$data[43,0] = "06:52:23"
$data[43,1] = "07:23"
$data[43,2] = "10_OLMOS"
$data[43,3] = 31
$data[43,4] = "LP1912"
$data[44,0] = "06:52:23"
$data[44,1] = "07:25"
$data[44,2] = "10_OLMOS"
$data[44,3] = 33
$data[44,4] = "LP1912"
$data[45,0] = "05:51:32"
$data[45,1] = "07:28"
$data[45,2] = "10_OLMOS"
$data[45,3] = 97
$data[45,4] = "LP1912"
$data[46,0] = "05:51:32"
$data[46,1] = "07:31"
$data[46,2] = "11_ETCHEVERRY"
$data[46,3] = 100
$data[46,4] = "LP1912"
$data[47,0] = "05:51:32"
$data[47,1] = "07:31"
$data[47,2] = "16_SANTA ANA"
$data[47,3] = 100
$data[47,4] = "LP1912"
$data[48,0] = "05:51:32"
$data[48,1] = "07:32"
$data[48,2] = "84_COLONIA URQUIZA-ESC 49"
$data[48,3] = 101
$data[48,4] = "LP1912"
$data[49,0] = "05:51:32"
$data[49,1] = "07:36"
$data[49,2] = "27_EL RETIRO"
$data[49,3] = 105
$data[49,4] = "LP1912"
$data[50,0] = "06:52:23"
$data[50,1] = "07:37"
$data[50,2] = "27_EL RETIRO"
$data[50,3] = 45
$data[50,4] = "LP1912"
$data[51,0] = "06:23:52"
$data[51,1] = "07:38"
$data[51,2] = "10_OLMOS"
$data[51,3] = 75
$data[51,4] = "LP1912"
$data[52,0] = "05:51:32"
$data[52,1] = "07:39"
$data[52,2] = "10_OLMOS"
$data[52,3] = 108
$data[52,4] = "LP1912"
$data[53,0] = "05:51:32"
$data[53,1] = "07:47"
$data[53,2] = "14_ABASTO"
$data[53,3] = 116
$data[53,4] = "LP1912"
$data[54,0] = "07:46:15"
$data[54,1] = "07:50"
$data[54,2] = "14_ABASTO"
$data[54,3] = 4
$data[54,4] = "LP1912"
$data[55,0] = "06:23:52"
$data[55,1] = "07:51"
$data[55,2] = "215D_EL PATO"
$data[55,3] = 88
$data[55,4] = "LP1912"
$data[56,0] = "07:46:15"
$data[56,1] = "07:55"
$data[56,2] = "10_OLMOS"
$data[56,3] = 9
$data[56,4] = "LP1912"
$data[57,0] = "07:46:15"
$data[57,1] = "07:56"
$data[57,2] = "16_SANTA ANA"
$data[57,3] = 10
$data[57,4] = "LP1912"
$data[58,0] = "07:59:05"
$data[58,1] = "08:02"
$data[58,2] = "16_SANTA ANA"
$data[58,3] = 3
$data[58,4] = "LP1912"
$data[59,0] = "06:52:23"
$data[59,1] = "08:03"
$data[59,2] = "23_HERNANDEZ"
$data[59,3] = 71
$data[59,4] = "LP1912"
$data[60,0] = "06:23:52"
$data[60,1] = "08:05"
$data[60,2] = "23_HERNANDEZ"
$data[60,3] = 102
$data[60,4] = "LP1912"
$data[61,0] = "07:46:15"
$data[61,1] = "08:09"
$data[61,2] = "11_ETCHEVERRY"
$data[61,3] = 23
$data[61,4] = "LP1912"
$data[62,0] = "07:59:05"
$data[62,1] = "08:11"
$data[62,2] = "11_ETCHEVERRY"
$data[62,3] = 12
$data[62,4] = "LP1912"
$data[63,0] = "06:23:52"
$data[63,1] = "08:12"
$data[63,2] = "15_ABASTO"
$data[63,3] = 109
$data[63,4] = "LP1912"
$data[64,0] = "06:23:52"
$data[64,1] = "08:20"
$data[64,2] = "26_HERNANDEZ"
$data[64,3] = 117
$data[64,4] = "LP1912"
$data[65,0] = "06:52:23"
$data[65,1] = "08:21"
$data[65,2] = "26_HERNANDEZ"
$data[65,3] = 89
$data[65,4] = "LP1912"
$data[66,0] = "06:23:52"
$data[66,1] = "08:22"
$data[66,2] = "16_P MOR-SANTA ANA"
$data[66,3] = 119
$data[66,4] = "LP1912"
$data[67,0] = "07:46:15"
$data[67,1] = "08:23"
$data[67,2] = "16_P MOR-SANTA ANA"
$data[67,3] = 37
$data[67,4] = "LP1912"
$data[68,0] = "06:52:23"
$data[68,1] = "08:23"
$data[68,2] = "215B_EL PATO"
$data[68,3] = 91
$data[68,4] = "LP1912"
$data[69,0] = "06:52:23"
$data[69,1] = "08:27"
$data[69,2] = "84_COLONIA URQUIZA-ESC 49"
$data[69,3] = 95
$data[69,4] = "LP1912"
$data[70,0] = "07:46:15"
$data[70,1] = "08:33"
$data[70,2] = "10_OLMOS"
$data[70,3] = 47
$data[70,4] = "LP1912"
$data[71,0] = "07:46:15"
$data[71,1] = "08:33"
$data[71,2] = "23_HERNANDEZ"
$data[71,3] = 47
$data[71,4] = "LP1912"
$data[72,0] = "07:46:15"
$data[72,1] = "08:34"
$data[72,2] = "26_HERNANDEZ"
$data[72,3] = 48
$data[72,4] = "LP1912"
$data[73,0] = "08:21:27"
$data[73,1] = "08:34"
$data[73,2] = "23_HERNANDEZ"
$data[73,3] = 13
$data[73,4] = "LP1912"
$data[74,0] = "07:59:05"
$data[74,1] = "08:35"
$data[74,2] = "23_HERNANDEZ"
$data[74,3] = 36
$data[74,4] = "LP1912"
$data[75,0] = "08:39:56"
$data[75,1] = "08:41"
$data[75,2] = "81_EL PELIGRO"
$data[75,3] = 2
$data[75,4] = "LP1912"
$data[76,0] = "08:39:56"
$data[76,1] = "08:42"
$data[76,2] = "14_ABASTO"
$data[76,3] = 3
$data[76,4] = "LP1912"
$data[77,0] = "06:52:23"
$data[77,1] = "08:42"
$data[77,2] = "81_EL PELIGRO"
$data[77,3] = 110
$data[77,4] = "LP1912"
$data[78,0] = "08:21:27"
$data[78,1] = "08:43"
$data[78,2] = "14_ABASTO"
$data[78,3] = 22
$data[78,4] = "LP1912"
$data[79,0] = "07:46:15"
$data[79,1] = "08:44"
$data[79,2] = "14_ABASTO"
$data[79,3] = 58
$data[79,4] = "LP1912"
$data[80,0] = "07:59:05"
$data[80,1] = "08:48"
$data[80,2] = "26_HERNANDEZ"
$data[80,3] = 49
$data[80,4] = "LP1912"
$data[81,0] = "08:39:56"
$data[81,1] = "08:52"
$data[81,2] = "10_OLMOS"
$data[81,3] = 13
$data[81,4] = "LP1912"
$data[82,0] = "07:59:05"
$data[82,1] = "08:53"
$data[82,2] = "10_OLMOS"
$data[82,3] = 54
$data[82,4] = "LP1912"
$data[83,0] = "07:46:15"
$data[83,1] = "08:54"
$data[83,2] = "17_ROMERO"
$data[83,3] = 68
$data[83,4] = "LP1912"
$data[84,0] = "08:21:27"
$data[84,1] = "09:01"
$data[84,2] = "23_HERNANDEZ"
$data[84,3] = 40
$data[84,4] = "LP1912"
$data[85,0] = "08:21:27"
$data[85,1] = "09:01"
$data[85,2] = "215A_EL PATO"
$data[85,3] = 40
$data[85,4] = "LP1912"
$data[86,0] = "07:46:15"
$data[86,1] = "09:02"
$data[86,2] = "215A_EL PATO"
$data[86,3] = 76
$data[86,4] = "LP1912"
$data[87,0] = "08:21:27"
$data[87,1] = "09:03"
$data[87,2] = "11_ETCHEVERRY"
$data[87,3] = 42
$data[87,4] = "LP1912"
$data[88,0] = "07:46:15"
$data[88,1] = "09:04"
$data[88,2] = "11_ETCHEVERRY"
$data[88,3] = 78
$data[88,4] = "LP1912"
$data[89,0] = "08:39:56"
$data[89,1] = "09:05"
$data[89,2] = "23_HERNANDEZ"
$data[89,3] = 26
$data[89,4] = "LP1912"
$data[90,0] = "08:21:27"
$data[90,1] = "09:09"
$data[90,2] = "26_HERNANDEZ"
$data[90,3] = 48
$data[90,4] = "LP1912"
$data[91,0] = "08:21:27"
$data[91,1] = "09:10"
$data[91,2] = "16_P MOR-SANTA ANA"
$data[91,3] = 49
$data[91,4] = "LP1912"
$data[92,0] = "07:46:15"
$data[92,1] = "09:11"
$data[92,2] = "16_P MOR-SANTA ANA"
$data[92,3] = 85
$data[92,4] = "LP1912"
$data[93,0] = "08:39:56"
$data[93,1] = "09:12"
$data[93,2] = "10_OLMOS"
$data[93,3] = 33
$data[93,4] = "LP1912"
$data[94,0] = "08:21:27"
$data[94,1] = "09:16"
$data[94,2] = "27_EL RETIRO"
$data[94,3] = 55
$data[94,4] = "LP1912"
$data[95,0] = "07:46:15"
$data[95,1] = "09:17"
$data[95,2] = "27_EL RETIRO"
$data[95,3] = 91
$data[95,4] = "LP1912"
$data[96,0] = "08:39:56"
$data[96,1] = "09:20"
$data[96,2] = "26_HERNANDEZ"
$data[96,3] = 41
$data[96,4] = "LP1912"
$data[97,0] = "07:46:15"
$data[97,1] = "09:21"
$data[97,2] = "26_HERNANDEZ"
$data[97,3] = 95
$data[97,4] = "LP1912"
$data[98,0] = "07:46:15"
$data[98,1] = "09:22"
$data[98,2] = "16_SANTA ANA"
$data[98,3] = 96
$data[98,4] = "LP1912"
$data[99,0] = "08:21:27"
$data[99,1] = "09:22"
$data[99,2] = "17_ROMERO"
$data[99,3] = 61
$data[99,4] = "LP1912"
$data[100,0] = "07:59:05"
$data[100,1] = "09:23"
$data[100,2] = "16_SANTA ANA"
$data[100,3] = 84
$data[100,4] = "LP1912"
$data[101,0] = "07:46:15"
$data[101,1] = "09:23"
$data[101,2] = "17_ROMERO"
$data[101,3] = 97
$data[101,4] = "LP1912"
$data[102,0] = "08:21:27"
$data[102,1] = "09:23"
$data[102,2] = "11_ETCHEVERRY"
$data[102,3] = 62
$data[102,4] = "LP1912"
$data[103,0] = "07:46:15"
$data[103,1] = "09:24"
$data[103,2] = "11_ETCHEVERRY"
$data[103,3] = 98
$data[103,4] = "LP1912"
$data[104,0] = "08:39:56"
$data[104,1] = "09:27"
$data[104,2] = "26_HERNANDEZ"
$data[104,3] = 48
$data[104,4] = "LP1912"
$data[105,0] = "07:46:15"
$data[105,1] = "09:32"
$data[105,2] = "15_ABASTO"
$data[105,3] = 106
$data[105,4] = "LP1912"
$data[106,0] = "07:46:15"
$data[106,1] = "09:33"
$data[106,2] = "10_OLMOS"
$data[106,3] = 107
$data[106,4] = "LP1912"
$data[107,0] = "08:39:56"
$data[107,1] = "09:34"
$data[107,2] = "23_HERNANDEZ"
$data[107,3] = 55
$data[107,4] = "LP1912"
$data[108,0] = "08:39:56"
$data[108,1] = "09:34"
$data[108,2] = "16_SANTA ANA"
$data[108,3] = 55
$data[108,4] = "LP1912"
$data[109,0] = "08:21:27"
$data[109,1] = "09:36"
$data[109,2] = "16_SANTA ANA"
$data[109,3] = 75
$data[109,4] = "LP1912"
$data[110,0] = "08:39:56"
$data[110,1] = "09:41"
$data[110,2] = "215C_EL PATO"
$data[110,3] = 62
$data[110,4] = "LP1912"
$data[111,0] = "07:46:15"
$data[111,1] = "09:42"
$data[111,2] = "215C_EL PATO"
$data[111,3] = 116
$data[111,4] = "LP1912"
$data[112,0] = "08:21:27"
$data[112,1] = "09:43"
$data[112,2] = "14_ABASTO"
$data[112,3] = 82
$data[112,4] = "LP1912"
$data[113,0] = "07:59:05"
$data[113,1] = "09:44"
$data[113,2] = "14_ABASTO"
$data[113,3] = 105
$data[113,4] = "LP1912"
$data[114,0] = "07:59:05"
$data[114,1] = "09:52"
$data[114,2] = "15_ABASTO"
$data[114,3] = 113
$data[114,4] = "LP1912"
$data[115,0] = "08:21:27"
$data[115,1] = "10:12"
$data[115,2] = "15_ABASTO"
$data[115,3] = 111
$data[115,4] = "LP1912"
$data[116,0] = "08:39:56"
$data[116,1] = "10:22"
$data[116,2] = "17_ROMERO"
$data[116,3] = 103
$data[116,4] = "LP1912"
$data[117,0] = "08:39:56"
$data[117,1] = "10:26"
$data[117,2] = "215A_EL PATO"
$data[117,3] = 107
$data[117,4] = "LP1912"
$ws.Range("A6:E123").Value = $data

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A2").Value = "Última actualización: 08:39:56"
$ws.Range("A3").Value = "Total filas: 18"
$data = New-Object 'object[,]' 18,5
$data[0,0] = "04:36:47"
$data[0,1] = "04:45"
$data[0,2] = "215A_EL PATO"
$data[0,3] = 9
$data[0,4] = "LP1912"
$data[1,0] = "04:01:01"
$data[1,1] = "04:46"
$data[1,2] = "215A_EL PATO"
$data[1,3] = 45
$data[1,4] = "LP1912"
$data[2,0] = "04:36:47"
$data[2,1] = "05:34"
$data[2,2] = "215B_EL PATO"
$data[2,3] = 58
$data[2,4] = "LP1912"
$data[3,0] = "04:01:01"
$data[3,1] = "05:35"
$data[3,2] = "215B_EL PATO"
$data[3,3] = 94
$data[3,4] = "LP1912"
$data[4,0] = "04:36:47"
$data[4,1] = "06:11"
$data[4,2] = "215A_EL PATO"
$data[4,3] = 95
$data[4,4] = "LP1912"
$data[5,0] = "05:20:00"
$data[5,1] = "06:12"
$data[5,2] = "215A_EL PATO"
$data[5,3] = 52
$data[5,4] = "LP1912"
$data[6,0] = "04:51:28"
$data[6,1] = "06:46"
$data[6,2] = "215C_EL PATO"
$data[6,3] = 115
$data[6,4] = "LP1912"
$data[7,0] = "05:20:00"
$data[7,1] = "06:47"
$data[7,2] = "215C_EL PATO"
$data[7,3] = 87
$data[7,4] = "LP1912"
$data[8,0] = "06:52:23"
$data[8,1] = "06:52"
$data[8,2] = "215C_EL PATO"
$data[8,3] = 0
$data[8,4] = "LP1912"
$data[9,0] = "05:51:32"
$data[9,1] = "07:11"
$data[9,2] = "215A_EL PATO"
$data[9,3] = 80
$data[9,4] = "LP1912"
$data[10,0] = "05:20:00"
$data[10,1] = "07:12"
$data[10,2] = "215A_EL PATO"
$data[10,3] = 112
$data[10,4] = "LP1912"
$data[11,0] = "06:23:52"
$data[11,1] = "07:51"
$data[11,2] = "215D_EL PATO"
$data[11,3] = 88
$data[11,4] = "LP1912"
$data[12,0] = "06:52:23"
$data[12,1] = "08:23"
$data[12,2] = "215B_EL PATO"
$data[12,3] = 91
$data[12,4] = "LP1912"
$data[13,0] = "08:21:27"
$data[13,1] = "09:01"
$data[13,2] = "215A_EL PATO"
$data[13,3] = 40
$data[13,4] = "LP1912"
$data[14,0] = "07:46:15"
$data[14,1] = "09:02"
$data[14,2] = "215A_EL PATO"
$data[14,3] = 76
$data[14,4] = "LP1912"
$data[15,0] = "08:39:56"
$data[15,1] = "09:41"
$data[15,2] = "215C_EL PATO"
$data[15,3] = 62
$data[15,4] = "LP1912"
$data[16,0] = "07:46:15"
$data[16,1] = "09:42"
$data[16,2] = "215C_EL PATO"
$data[16,3] = 116
$data[16,4] = "LP1912"
$data[17,0] = "08:39:56"
$data[17,1] = "10:26"
$data[17,2] = "215A_EL PATO"
$data[17,3] = 107
$data[17,4] = "LP1912"
$ws.Range("A6:E23").Value = $data

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 08:39:56"
$ws.Range("A3").Value = "Total filas: 26"
$data = New-Object 'object[,]' 26,5
$data[0,0] = "04:36:47"
$data[0,1] = "05:43"
$data[0,2] = "215A_LA PLATA"
$data[0,3] = 67
$data[0,4] = "L6173"
$data[1,0] = "04:01:01"
$data[1,1] = "05:44"
$data[1,2] = "215A_LA PLATA"
$data[1,3] = 103
$data[1,4] = "L6173"
$data[2,0] = "04:51:28"
$data[2,1] = "06:09"
$data[2,2] = "215A_LA PLATA"
$data[2,3] = 78
$data[2,4] = "L6173"
$data[3,0] = "04:36:47"
$data[3,1] = "06:10"
$data[3,2] = "215A_LA PLATA"
$data[3,3] = 94
$data[3,4] = "L6173"
$data[4,0] = "04:36:47"
$data[4,1] = "06:32"
$data[4,2] = "215C_LA PLATA"
$data[4,3] = 116
$data[4,4] = "L6203"
$data[5,0] = "04:51:28"
$data[5,1] = "06:33"
$data[5,2] = "215C_LA PLATA"
$data[5,3] = 102
$data[5,4] = "L6203"
$data[6,0] = "06:23:52"
$data[6,1] = "06:59"
$data[6,2] = "215B_LP-P MOR-1 Y 57"
$data[6,3] = 36
$data[6,4] = "L6173"
$data[7,0] = "05:20:00"
$data[7,1] = "07:00"
$data[7,2] = "215B_LP-P MOR-1 Y 57"
$data[7,3] = 100
$data[7,4] = "L6173"
$data[8,0] = "05:51:32"
$data[8,1] = "07:35"
$data[8,2] = "215A_LA PLATA"
$data[8,3] = 104
$data[8,4] = "L6173"
$data[9,0] = "06:52:23"
$data[9,1] = "07:38"
$data[9,2] = "215A_LA PLATA"
$data[9,3] = 46
$data[9,4] = "L6173"
$data[10,0] = "07:16:53"
$data[10,1] = "07:44"
$data[10,2] = "215A_LA PLATA"
$data[10,3] = 28
$data[10,4] = "L6173"
$data[11,0] = "07:46:15"
$data[11,1] = "07:51"
$data[11,2] = "215A_LA PLATA"
$data[11,3] = 5
$data[11,4] = "L6173"
$data[12,0] = "06:23:52"
$data[12,1] = "08:06"
$data[12,2] = "215C_LA PLATA"
$data[12,3] = 103
$data[12,4] = "L6203"
$data[13,0] = "07:46:15"
$data[13,1] = "08:09"
$data[13,2] = "215C_LA PLATA"
$data[13,3] = 23
$data[13,4] = "L6203"
$data[14,0] = "07:16:53"
$data[14,1] = "08:10"
$data[14,2] = "215C_LA PLATA"
$data[14,3] = 54
$data[14,4] = "L6203"
$data[15,0] = "06:52:23"
$data[15,1] = "08:11"
$data[15,2] = "215C_LA PLATA"
$data[15,3] = 79
$data[15,4] = "L6203"
$data[16,0] = "07:59:05"
$data[16,1] = "08:15"
$data[16,2] = "215C_LA PLATA"
$data[16,3] = 16
$data[16,4] = "L6203"
$data[17,0] = "06:52:23"
$data[17,1] = "08:40"
$data[17,2] = "215A_LA PLATA"
$data[17,3] = 108
$data[17,4] = "L6173"
$data[18,0] = "07:46:15"
$data[18,1] = "08:45"
$data[18,2] = "215A_LA PLATA"
$data[18,3] = 59
$data[18,4] = "L6173"
$data[19,0] = "07:59:05"
$data[19,1] = "08:46"
$data[19,2] = "215A_LA PLATA"
$data[19,3] = 47
$data[19,4] = "L6173"
$data[20,0] = "08:21:27"
$data[20,1] = "08:52"
$data[20,2] = "215A_LA PLATA"
$data[20,3] = 31
$data[20,4] = "L6173"
$data[21,0] = "08:39:56"
$data[21,1] = "08:55"
$data[21,2] = "215A_LA PLATA"
$data[21,3] = 16
$data[21,4] = "L6173"
$data[22,0] = "07:16:53"
$data[22,1] = "09:08"
$data[22,2] = "215D_LA PLATA"
$data[22,3] = 112
$data[22,4] = "L6203"
$data[23,0] = "07:46:15"
$data[23,1] = "09:09"
$data[23,2] = "215D_LA PLATA"
$data[23,3] = 83
$data[23,4] = "L6203"
$data[24,0] = "08:39:56"
$data[24,1] = "10:02"
$data[24,2] = "215B_LP-P MOR-40 Y 115"
$data[24,3] = 83
$data[24,4] = "L6173"
$data[25,0] = "08:21:27"
$data[25,1] = "10:03"
$data[25,2] = "215B_LP-P MOR-40 Y 115"
$data[25,3] = 102
$data[25,4] = "L6173"
$ws.Range("A6:E31").Value = $data

